$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # ALC
$ws2 = $wb.Worksheets.Item(2)  # ARM
$ws3 = $wb.Worksheets.Item(3)  # BSM
$ws4 = $wb.Worksheets.Item(4)  # CRP
$ws5 = $wb.Worksheets.Item(5)  # CUL
$ws6 = $wb.Worksheets.Item(6)  # GSM
$ws7 = $wb.Worksheets.Item(7)  # LTW
$ws8 = $wb.Worksheets.Item(8)  # WVR

# ALC row 11: Gotta Bounce | Rubber
$ws1.Cells.Item(11, 8).Value = 36996.9
$ws1.Cells.Item(11, 9).Value = 36996.9
$ws1.Cells.Item(11, 11).Value = 36996.9
$ws1.Cells.Item(11, 13).Value = -36856.9

# ALC row 43: Growing Is Knowing | Growth Formula Gamma
$ws1.Cells.Item(43, 8).Value = 2735.2
$ws1.Cells.Item(43, 9).Value = 2575.3333
$ws1.Cells.Item(43, 11).Value = 2575.3333
$ws1.Cells.Item(43, 13).Value = -2506.3333

# ALC row 55: A Real Smooth Move | Lanolin
$ws1.Cells.Item(55, 8).Value = 195
$ws1.Cells.Item(55, 9).Value = 195
$ws1.Cells.Item(55, 10).Value = 0
$ws1.Cells.Item(55, 11).Value = 195
$ws1.Cells.Item(55, 12).Value = 0
$ws1.Cells.Item(55, 13).Value = 19
$ws1.Cells.Item(55, 14).ClearContents()

# ALC row 88: The Grave of Hemlock Groves | Growth Formula Zeta
$ws1.Cells.Item(88, 8).Value = 1660.2142
$ws1.Cells.Item(88, 9).Value = 2796
$ws1.Cells.Item(88, 10).Value = 1029.2222
$ws1.Cells.Item(88, 11).Value = 2796
$ws1.Cells.Item(88, 12).Value = 1029.2222
$ws1.Cells.Item(88, 13).Value = -2390
$ws1.Cells.Item(88, 14).Value = -1841.2222

# ALC row 91: Dappling the Highlands (L) | Growth Formula Zeta
$ws1.Cells.Item(91, 8).Value = 1660.2142
$ws1.Cells.Item(91, 9).Value = 2796
$ws1.Cells.Item(91, 10).Value = 1029.2222
$ws1.Cells.Item(91, 11).Value = 2796
$ws1.Cells.Item(91, 12).Value = 1029.2222
$ws1.Cells.Item(91, 13).Value = -1392
$ws1.Cells.Item(91, 14).Value = -3837.2222

# ALC row 98: The Dotted Line | Enchanted Durium Ink
$ws1.Cells.Item(98, 8).Value = 772.8125
$ws1.Cells.Item(98, 9).Value = 772.8125
$ws1.Cells.Item(98, 11).Value = 772.8125
$ws1.Cells.Item(98, 13).Value = 725.1875

# ALC row 121: Mindful Medicine | Tincture of Mind
$ws1.Cells.Item(121, 8).Value = 3663.2307
$ws1.Cells.Item(121, 10).Value = 3663.2307
$ws1.Cells.Item(121, 12).Value = 10989.6921
$ws1.Cells.Item(121, 14).Value = -14483.6921

# ALC row 122: Wishful Inking | Enchanted High Durium Ink
$ws1.Cells.Item(122, 8).Value = 772.8125
$ws1.Cells.Item(122, 9).Value = 772.8125
$ws1.Cells.Item(122, 11).Value = 2318.4375
$ws1.Cells.Item(122, 13).Value = 131.5625

# ALC row 125: Body over Mind | Grade 5 Dexterity Alkahest
$ws1.Cells.Item(125, 8).Value = 1499.5
$ws1.Cells.Item(125, 9).Value = 1499
$ws1.Cells.Item(125, 11).Value = 13491
$ws1.Cells.Item(125, 13).Value = -11031

# ALC row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws1.Cells.Item(137, 8).Value = 6897666
$ws1.Cells.Item(137, 9).Value = 7143936
$ws1.Cells.Item(137, 11).Value = 21431808
$ws1.Cells.Item(137, 13).Value = -21429258

# ALC row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws1.Cells.Item(138, 8).Value = 5052.878
$ws1.Cells.Item(138, 9).Value = 2710.625
$ws1.Cells.Item(138, 10).Value = 5620.697
$ws1.Cells.Item(138, 11).Value = 8131.875
$ws1.Cells.Item(138, 12).Value = 16862.091
$ws1.Cells.Item(138, 13).Value = -2991.875
$ws1.Cells.Item(138, 14).Value = -27142.091

# ALC row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws1.Cells.Item(141, 8).Value = 6051.294
$ws1.Cells.Item(141, 9).Value = 2322.6667
$ws1.Cells.Item(141, 11).Value = 6968.000100000001
$ws1.Cells.Item(141, 13).Value = -1788.000100000001

# ARM row 110: Scheduled Maintenance | Deepgold Ingot
$ws2.Cells.Item(110, 8).Value = 1196.3334
$ws2.Cells.Item(110, 9).Value = 849.5
$ws2.Cells.Item(110, 10).Value = 1890
$ws2.Cells.Item(110, 11).Value = 849.5
$ws2.Cells.Item(110, 12).Value = 1890
$ws2.Cells.Item(110, 13).Value = 1195.5
$ws2.Cells.Item(110, 14).Value = -5980

# ARM row 111: Hedging Bets | Deepgold Surcoat of Maiming
$ws2.Cells.Item(111, 8).Value = 63332.668
$ws2.Cells.Item(111, 10).Value = 63332.668
$ws2.Cells.Item(111, 12).Value = 63332.668
$ws2.Cells.Item(111, 14).Value = -71512.66800000001

# BSM row 86: Through Thick and Thin | Adamantite Nugget
$ws3.Cells.Item(86, 8).Value = 2438.2856
$ws3.Cells.Item(86, 9).Value = 1976.091
$ws3.Cells.Item(86, 11).Value = 1976.091
$ws3.Cells.Item(86, 13).Value = -853.0909999999999

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws3.Cells.Item(89, 8).Value = 2438.2856
$ws3.Cells.Item(89, 9).Value = 1976.091
$ws3.Cells.Item(89, 11).Value = 9880.455
$ws3.Cells.Item(89, 13).Value = -4264.455

# BSM row 99: Meddle in Metal | Oroshigane Ingot
$ws3.Cells.Item(99, 8).Value = 9305.529
$ws3.Cells.Item(99, 9).Value = 10514.154
$ws3.Cells.Item(99, 11).Value = 10514.154
$ws3.Cells.Item(99, 13).Value = -9016.154

# BSM row 107: The Gold Experience | Deepgold Nugget
$ws3.Cells.Item(107, 8).Value = 1420.5333
$ws3.Cells.Item(107, 9).Value = 1524.5834
$ws3.Cells.Item(107, 10).Value = 1004.3333
$ws3.Cells.Item(107, 11).Value = 1524.5834
$ws3.Cells.Item(107, 12).Value = 1004.3333
$ws3.Cells.Item(107, 13).Value = 395.4166
$ws3.Cells.Item(107, 14).Value = -4844.3333

# BSM row 119: Bae Blade | Dwarven Mythril Uchigatana
$ws3.Cells.Item(119, 8).Value = 98333.336
$ws3.Cells.Item(119, 10).Value = 98333.336
$ws3.Cells.Item(119, 12).Value = 98333.336
$ws3.Cells.Item(119, 14).Value = -108009.336

# BSM row 134: Ruthenium Supremium | Ruthenium Ingot
$ws3.Cells.Item(134, 8).Value = 4730769.5
$ws3.Cells.Item(134, 9).Value = 3970826.5
$ws3.Cells.Item(134, 11).Value = 11912479.5
$ws3.Cells.Item(134, 13).Value = -11909944.5

# CRP row 4: A Clogful of Camaraderie | Maple Clogs
$ws4.Cells.Item(4, 8).Value = 866.3333
$ws4.Cells.Item(4, 10).Value = 866.3333
$ws4.Cells.Item(4, 12).Value = 866.3333
$ws4.Cells.Item(4, 14).Value = -1090.3333

# CRP row 31: Wall Not Found | Walnut Lumber
$ws4.Cells.Item(31, 8).Value = 1010486.75
$ws4.Cells.Item(31, 10).Value = 6344.75
$ws4.Cells.Item(31, 12).Value = 6344.75
$ws4.Cells.Item(31, 14).Value = -6934.75

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws4.Cells.Item(34, 8).Value = 1010486.75
$ws4.Cells.Item(34, 10).Value = 6344.75
$ws4.Cells.Item(34, 12).Value = 6344.75
$ws4.Cells.Item(34, 14).Value = -6748.75

# CRP row 105: Zelkova, My Love | Zelkova Lumber
$ws4.Cells.Item(105, 8).Value = 3766.85
$ws4.Cells.Item(105, 9).Value = 4130.1113
$ws4.Cells.Item(105, 11).Value = 4130.1113
$ws4.Cells.Item(105, 13).Value = -2383.1113

# CRP row 112: Understaffed | Applewood Cane
$ws4.Cells.Item(112, 8).Value = 0
$ws4.Cells.Item(112, 10).Value = 0
$ws4.Cells.Item(112, 12).Value = 0
$ws4.Cells.Item(112, 14).ClearContents()

# CRP row 134: Wood You Be Quiet | Ceiba Lumber
$ws4.Cells.Item(134, 8).Value = 5015.75
$ws4.Cells.Item(134, 9).Value = 3556.4211
$ws4.Cells.Item(134, 11).Value = 10669.2633
$ws4.Cells.Item(134, 13).Value = -8134.263300000001

# CUL row 12: Butter Me Up | Kukuru Butter
$ws5.Cells.Item(12, 8).Value = 668.0833
$ws5.Cells.Item(12, 9).Value = 1682
$ws5.Cells.Item(12, 10).Value = 330.1111
$ws5.Cells.Item(12, 11).Value = 5046
$ws5.Cells.Item(12, 12).Value = 990.3333
$ws5.Cells.Item(12, 13).Value = -4873
$ws5.Cells.Item(12, 14).Value = -1336.3333

# CUL row 81: It Goes Down Smoothly | Frozen Spirits
$ws5.Cells.Item(81, 8).Value = 2503
$ws5.Cells.Item(81, 9).Value = 997
$ws5.Cells.Item(81, 10).Value = 3005
$ws5.Cells.Item(81, 11).Value = 2991
$ws5.Cells.Item(81, 12).Value = 9015
$ws5.Cells.Item(81, 13).Value = -1868
$ws5.Cells.Item(81, 14).Value = -11261

# CUL row 84: Quenching the Flame (L) | Frozen Spirits
$ws5.Cells.Item(84, 8).Value = 2503
$ws5.Cells.Item(84, 9).Value = 997
$ws5.Cells.Item(84, 10).Value = 3005
$ws5.Cells.Item(84, 11).Value = 8973
$ws5.Cells.Item(84, 12).Value = 27045
$ws5.Cells.Item(84, 13).Value = -3357
$ws5.Cells.Item(84, 14).Value = -38277

# CUL row 93: Loquacious | Loquat Juice
$ws5.Cells.Item(93, 8).Value = 2249.8333
$ws5.Cells.Item(93, 9).Value = 2099
$ws5.Cells.Item(93, 11).Value = 6297
$ws5.Cells.Item(93, 13).Value = -4425

# CUL row 107: Slippery Service | Frantoio Oil
$ws5.Cells.Item(107, 8).Value = 4435.077
$ws5.Cells.Item(107, 9).Value = 329
$ws5.Cells.Item(107, 10).Value = 5412.7144
$ws5.Cells.Item(107, 11).Value = 987
$ws5.Cells.Item(107, 12).Value = 16238.1432
$ws5.Cells.Item(107, 13).Value = 933
$ws5.Cells.Item(107, 14).Value = -20078.1432

# CUL row 109: Cure for What Ails | Purple Carrot Juice
$ws5.Cells.Item(109, 8).Value = 5635.353
$ws5.Cells.Item(109, 9).Value = 2106.75
$ws5.Cells.Item(109, 10).Value = 8771.888999999999
$ws5.Cells.Item(109, 11).Value = 6320.25
$ws5.Cells.Item(109, 12).Value = 26315.667
$ws5.Cells.Item(109, 13).Value = -5280.25
$ws5.Cells.Item(109, 14).Value = -28395.667

# GSM row 20: Brothers in Arms | Brass Wristlets of Crafting
$ws6.Cells.Item(20, 8).Value = 4924.1577
$ws6.Cells.Item(20, 9).Value = 4197.2666
$ws6.Cells.Item(20, 10).Value = 7650
$ws6.Cells.Item(20, 11).Value = 4197.2666
$ws6.Cells.Item(20, 12).Value = 7650
$ws6.Cells.Item(20, 13).Value = -3952.2666
$ws6.Cells.Item(20, 14).Value = -8140

# GSM row 80: Needs More Prayerbell | Hardsilver Ingot
$ws6.Cells.Item(80, 8).Value = 6840.25
$ws6.Cells.Item(80, 10).Value = 6573.3335
$ws6.Cells.Item(80, 12).Value = 6573.3335
$ws6.Cells.Item(80, 14).Value = -8569.333500000001

# GSM row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws6.Cells.Item(83, 8).Value = 6840.25
$ws6.Cells.Item(83, 10).Value = 6573.3335
$ws6.Cells.Item(83, 12).Value = 32866.6675
$ws6.Cells.Item(83, 14).Value = -42850.6675

# GSM row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws6.Cells.Item(97, 8).Value = 1194.6666
$ws6.Cells.Item(97, 9).Value = 1084.1666
$ws6.Cells.Item(97, 10).Value = 1283.0667
$ws6.Cells.Item(97, 11).Value = 1084.1666
$ws6.Cells.Item(97, 12).Value = 1283.0667
$ws6.Cells.Item(97, 13).Value = -588.1666
$ws6.Cells.Item(97, 14).Value = -2275.0667

# GSM row 132: On Board for Lar | Lar Ingot
$ws6.Cells.Item(132, 8).Value = 12791.697
$ws6.Cells.Item(132, 9).Value = 7290.9644
$ws6.Cells.Item(132, 11).Value = 21872.8932
$ws6.Cells.Item(132, 13).Value = -19342.8932

# LTW row 40: Best Served Toad | Toad Leather
$ws7.Cells.Item(40, 8).Value = 2600.1892
$ws7.Cells.Item(40, 9).Value = 2605.75
$ws7.Cells.Item(40, 11).Value = 2605.75
$ws7.Cells.Item(40, 13).Value = -2469.75

# LTW row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws7.Cells.Item(55, 8).Value = 1785.8148
$ws7.Cells.Item(55, 9).Value = 819.3077
$ws7.Cells.Item(55, 10).Value = 2683.2856
$ws7.Cells.Item(55, 11).Value = 819.3077
$ws7.Cells.Item(55, 12).Value = 2683.2856
$ws7.Cells.Item(55, 13).Value = -646.3077
$ws7.Cells.Item(55, 14).Value = -3029.2856

# LTW row 61: Spelling Me Softly | Raptor Leather
$ws7.Cells.Item(61, 8).Value = 11717.786
$ws7.Cells.Item(61, 9).Value = 10115.2
$ws7.Cells.Item(61, 11).Value = 10115.2
$ws7.Cells.Item(61, 13).Value = -9913.200000000001

# LTW row 113: Peace in Rest | Atrociraptor Leather
$ws7.Cells.Item(113, 8).Value = 11717.786
$ws7.Cells.Item(113, 9).Value = 10115.2
$ws7.Cells.Item(113, 11).Value = 10115.2
$ws7.Cells.Item(113, 13).Value = -7945.200000000001

# LTW row 136: Respect for Br'aax | Br'aax Leather
$ws7.Cells.Item(136, 8).Value = 14585143
$ws7.Cells.Item(136, 9).Value = 14064080
$ws7.Cells.Item(136, 10).Value = 16669393
$ws7.Cells.Item(136, 11).Value = 42192240
$ws7.Cells.Item(136, 12).Value = 50008179
$ws7.Cells.Item(136, 13).Value = -42189690
$ws7.Cells.Item(136, 14).Value = -50013279

# LTW row 140: Worqor Zormor or Bust | Gargantuaskin Shoes of Healing
$ws7.Cells.Item(140, 8).Value = 99994
$ws7.Cells.Item(140, 10).Value = 99994
$ws7.Cells.Item(140, 12).Value = 99994
$ws7.Cells.Item(140, 14).Value = -110354

# WVR row 100: Of Great Import | Kudzu Thread
$ws8.Cells.Item(100, 8).Value = 815.4545000000001
$ws8.Cells.Item(100, 9).Value = 497
$ws8.Cells.Item(100, 10).Value = 4000
$ws8.Cells.Item(100, 11).Value = 994
$ws8.Cells.Item(100, 12).Value = 8000
$ws8.Cells.Item(100, 13).Value = -453
$ws8.Cells.Item(100, 14).Value = -9082

# WVR row 113: A Tender Table | Pixie Floss
$ws8.Cells.Item(113, 8).Value = 2035.9524
$ws8.Cells.Item(113, 9).Value = 1953.1875
$ws8.Cells.Item(113, 10).Value = 2300.8
$ws8.Cells.Item(113, 11).Value = 5859.5625
$ws8.Cells.Item(113, 12).Value = 6902.400000000001
$ws8.Cells.Item(113, 13).Value = -3689.5625
$ws8.Cells.Item(113, 14).Value = -11242.4

# WVR row 136: Weaving the Envelope | Sarcenet Cloth
$ws8.Cells.Item(136, 8).Value = 1767937
$ws8.Cells.Item(136, 9).Value = 1233054.1
$ws8.Cells.Item(136, 11).Value = 3699162.3
$ws8.Cells.Item(136, 13).Value = -3696612.3
